# Adds two new columns, I (I0) and J (IF), to the sheet.
# Header row (row 1) gets "I0" in I1 and "IF" in J1, styled like the
# other header cells (style index 1 -> bold font, border, centered).
# Rows 2-46 get the numeric values for the new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells, copy style from existing header cell (H1) so the new
# headers look consistent with the rest of row 1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2..46 (I column then J column)
$iValues = @(5,8,8,7,9,7,6,8,8,6,8,7,7,6,5,7,9,9,8,7,7,8,8,8,7,7,6,6,8,9,6,6,6,5,9,6,7,9,9,5,9,7,6,5,3)
$jValues = @(6,8,9,8,10,7,7,8,8,6,8,8,7,6,6,7,9,9,8,7,7,8,8,8,7,7,6,6,8,9,6,7,7,6,9,6,7,9,9,6,9,7,6,5,3)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
